# Add two new rows ("Dog" / "Snake") to the Translations sheet, each
# marked with the new "New" status, mirroring the existing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Dog"
$ws.Range("B6").Value = "Dog"
$ws.Range("C6").Value = "Dog"
$ws.Range("D6").Value = "New"

$ws.Range("A7").Value = "Snake"
$ws.Range("B7").Value = "Snake"
$ws.Range("C7").Value = "Snake"
$ws.Range("D7").Value = "New"
